$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# About sheet: A2 and A6 contain the build timestamp text
$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"
$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Black Eagle Coal Mine, United States, M3402, version ''mines - January 30 (built on ' + $newStamp + ')''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# Boundaries and methane sources sheet: column S (build_version) rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $cell = $wsBoundaries.Cells.Item($r, 19)
    $cell.Value = "mines - January 30 (built on " + $newStamp + ")"
}
